$wb = $excel.ActiveWorkbook

# --- Update "Invoice Details" sheet: revised confidence scores in H2:J2 ---
$wsDetails = $wb.Worksheets.Item("Invoice Details")
$wsDetails.Range("H2").Value = 0.78
$wsDetails.Range("I2").Value = 0.84
$wsDetails.Range("J2").Value = 0.96

# --- Populate "Line Items" sheet with header row + one extracted data row ---
$wsLine = $wb.Worksheets.Item("Line Items")

# Header row (bold, centered, bordered - matching the workbook's header style)
$wsLine.Range("A1").Value = "serial_number"
$wsLine.Range("B1").Value = "description"
$wsLine.Range("C1").Value = "quantity"
$wsLine.Range("D1").Value = "unit_price"
$wsLine.Range("E1").Value = "total_amount"
$wsLine.Range("F1").Value = "Confidence Score"
$wsLine.Range("G1").Value = "Validation Passed"

$headerRange = $wsLine.Range("A1:G1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160
$headerRange.Borders.LineStyle = 1

# Data row - numeric-looking OCR fields are kept as text, matching the source data
$wsLine.Range("A2").NumberFormat = "@"
$wsLine.Range("A2").Value = "1000"

$wsLine.Range("B2").Value = "YN 1000 4 1000`nI S`nYN 1000 1000`nMSH`nYN"

$wsLine.Range("C2").NumberFormat = "@"
$wsLine.Range("C2").Value = "1000"

$wsLine.Range("D2").NumberFormat = "@"
$wsLine.Range("D2").Value = "4"

$wsLine.Range("E2").NumberFormat = "@"
$wsLine.Range("E2").Value = "4000"

$wsLine.Range("F2").Value = 0.96
$wsLine.Range("G2").Value = $true
